$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append new row 43: 122. Best Time to Buy and Sell Stock II -------------
# Shared strings get appended in the order the distinct new text values are
# first written, so write the hyperlink URL (E) before the title (A) and the
# notes (D) to reproduce the original author's insertion order.
$url   = "https://leetcode.com/problems/best-time-to-buy-and-sell-stock-ii/solutions/208241/explanation-for-the-dummy-like-me/?envType=study-plan-v2&envId=top-interview-150 "
$title = "122. Best Time to Buy and Sell Stock II"
$notes = "The crux is that we solve this with a greedy approach. It is always better to buy local minima and sell at the next local high, rather than buy at the global minimum and sell at the global maximum. Use a while loop i < n, and 2 inner while loops for buy and sell, using prices[i+1] to search for conditions."

$ws.Range("E43").Value = $url
$ws.Range("A43").Value = $title
$ws.Range("B43").Value = "Medium"
$ws.Range("C43").Value = "Arrays"
$ws.Range("D43").Value = $notes

# Register the real hyperlink (adds the relationship + hyperlink entry), then
# restore the shared "Hyperlink" cell style used by the rest of column E.
$ws.Hyperlinks.Add($ws.Range("E43"), $url)
$ws.Range("E43").Style = "Hyperlink"

# Match the "Medium" difficulty shading used elsewhere in column B.
$ws.Range("B43").Interior.Color = 49407

# --- Update the saved view / selection ---------------------------------
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D49").Select()
